$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('B2').Value = 'PredicciA3n
                        para hoy                        15 de Agosto de 2023
Donde otras personas ven montanas que superar, tu ves soluciones. Cuando los demas estamos aturdidos y desconcertados, tu puedes encontrar el camino a traves de la confusion. Lo frustrante es que tus soluciones no siempre son bien recibidas. Por desgracia, solo porque tu lo veas claro, no significa que los otros esten dispuestos a seguirte. Asi que el hecho de que todo el mundo te pida consejo (y tu no sepas que hacer) es poco habitual, por no decir inquietante. Lo que era oscuro se aclarara bajo el influjo de la oscuridad de la Luna - y tu solucion te impresionara incluso a ti.'
$ws.Rows.Item(2).AutoFit()

$ws.Range('B3').Value = 'PredicciA3n
                        para hoy                        15 de Agosto de 2023
Lo malo de las nuevas ideas es que exigen un cambio de actitud. Ponerlas en practica significa reajustar las expectativas y alterar las rutinas. Lo cual esta bien cuando somos nosotros los promotores de la idea. Pero cuando viene de otra persona es mas dificil y tendemos a considerarla una molestia. Ademas tenemos problemas con las ideas que nacen de la necesidad (en lugar de la inspiracion). Sin embargo, venga de donde venga una idea, merece la pena investigarla bajo el influjo de la oscuridad de la Luna. Estate abierta a lo que venga. Una idea que descartarias facilmente podria ser revolucionaria.'
$ws.Rows.Item(3).AutoFit()

$ws.Range('B4').Value = 'PredicciA3n
                        para hoy                        15 de Agosto de 2023
Es como si estuvieras intentando montar un puzle sin tener la foto como guia. O como armar un mueble sin tener las instrucciones. Ni siquiera estas segura de si son las piezas correctas porque, hagas lo que hagas, las que tienes no parecen encajar entre si. !Que rabia! Mira que es frustrante. Pero tienes que seguir experimentando. No pierdas el optimismo. Puedes encontrar la respuesta a tu pregunta probando, equivocandote, mejorando y usando tu ingenio. Y una vez que la hayas encontrado, sabras que hacer exactamente para mejorar tus circunstancias actuales. !Sigue adelante!'
$ws.Rows.Item(4).AutoFit()

$ws.Range('B5').Value = 'PredicciA3n
                        para hoy                        15 de Agosto de 2023
La oscuridad de la Luna es un momento poderoso para ti. Como tu regente casi desaparece de los cielos, tus emociones se intensifican. Experimentas sentimientos profundos. Tus experiencias te parecen mas extremas. Lo que significa que las oportunidades te dan miedo. Pero estos sintomas son positivos. Te estan resultando muy saludables. Estas plantando cara a tus temores y venciendo tus inhibiciones. Asi que se valiente y centrate en esas dificultades. Estas al borde de un gran avance. Las repercusiones podrian suponer un cambio que sera beneficioso durante mucho tiempo.'
$ws.Rows.Item(5).AutoFit()

$ws.Range('B6').Value = 'PredicciA3n
                        para hoy                        15 de Agosto de 2023
Si tuvieras una maquina del tiempo, ?que harias? ?A que acontecimientos de tu pasado viajarias para cambiarlos? ?Que consejo le darias a tu antiguo yo? ?Tienes algun truco o idea que transmitirle? ?Algunas palabras de consuelo y aliento? A pesar de nuestros avances tecnicos, no se ha inventado ningun dispositivo de este tipo. Sin embargo, a medida que nos acercamos a la Luna Nueva (que sera en tu signo) algo con lo que nunca te has atrevido a sonar esta lo suficientemente cerca como para poder tocarlo. Algo imposible se esta haciendo posible. Tu perspectiva de futuro puede transformarse.'
$ws.Rows.Item(6).AutoFit()

$ws.Range('B7').Value = 'PredicciA3n
                        para hoy                        15 de Agosto de 2023
Si no te importa olvidar todo lo que te hayan contado acerca de la famosa capacidad para discernir que teneis los nativos de Virgo y vuestra minuciosidad, seria un signo perfecto. No es que esta imagen sea falsa... pero desde luego no esta completa. Y corres el peligro de que, al sentirte presionada a ajustarte a ciertas expectativas, puedas acabar quitando el brillo a tus otras cualidades, ocultandolas bajo una capa de modestia. La oscuridad de la Luna fomenta tu lado creativo y tu caracter desbordante de vida y entusiasmo. No te costaria nada dejar que las exigencias de los demas se interpusieran en tu camino. No lo permitas.'
$ws.Rows.Item(7).AutoFit()

$ws.Range('B8').Value = 'PredicciA3n
                        para hoy                        15 de Agosto de 2023
Lo que es importante para mi puede no serlo para ti. Que a mi me guste la pizza no significa que a ti tambien te guste. Lo que para mi puede ser una crisis, para ti puede ser un simple contratiempo. Sabemos esto. Pero se nos olvida. Por lo que sea, pensamos que todo el mundo siente y piensa igual que nosotros. Y luego nos preguntamos por que nadie nos comprende. Cierta persona de tu mundo esta tan interesada en que entiendas su punto de vista que hace todo lo posible por compartirlo contigo. Ten con ella toda la paciencia de que seas capaz. Significa mucho para ella. Hasta es posible acabe siendo importante para ti.'
$ws.Rows.Item(8).AutoFit()

$ws.Range('B9').Value = 'PredicciA3n
                        para hoy                        15 de Agosto de 2023
Si estas corriendo y ya llevas una pelota, va a ser todo un reto para ti coger otra pelota que casualmente rueda en tu direccion. ?Deberias soltar la primera? ?Como vas a saber a cual te conviene aferrarte: a la que ya tienes o a la que podrias tener si no hubieras cogido la otra antes? Vaya pregunta... es dificil dar una respuesta. Procura no confundir mas las cosas, dejando que te preocupe lo que piensen los demas. Ni tampoco dejes que te influya un equivocado sentido de la lealtad. Se justa. Pero no seas sentimental.'
$ws.Rows.Item(9).AutoFit()

$ws.Range('B10').Value = 'PredicciA3n
                        para hoy                        15 de Agosto de 2023
"La esperanza es esa cosa con plumas que se posa en el alma". Cito estas sabias palabras porque estas contemplando dejar aparcada una idea que podria estar a punto de alzar el vuelo. Nadie quiere sentirse decepcionado. Pero cerrarte a la oportunidad de poner en marcha un cambio, no cambiara nada. No puedes quedarte quieta y tampoco quieres hacerlo. Por supuesto que dar el primer paso te produce temor. Hara que te sientas vulnerable. Pero si eres lo bastante valiente como para desplegar tus alas, enseguida estaras volando mas alto de lo que te atreviste a sonar.'
$ws.Rows.Item(10).AutoFit()

$ws.Range('B11').Value = 'PredicciA3n
                        para hoy                        15 de Agosto de 2023
"?Juras decir toda la verdad y nada mas que la verdad?" !Esto si que es un reto! Puede parecer sencillo. Pero si tu y yo somos testigos del mismo suceso, cada uno le daremos, sin proponernoslo, nuestro propio giro a la historia. No es que queramos faltar a la verdad, pero no podemos evitar verlo a traves de nuestros ojos, lo que significa que nuestra experiencia (antes, durante e incluso despues del suceso) influye en nuestra version. Hoy ten esto en cuenta cuando intentes tomar una decision basandote en la version que ha dado otra persona de lo que ha pasado. No partas de ideas preconcebidas.'
$ws.Rows.Item(11).AutoFit()

$ws.Range('B12').Value = 'PredicciA3n
                        para hoy                        15 de Agosto de 2023
Es facil idealizar a la gente que nos importa. Pensar en su cuidado y apoyo puede emocionarnos. Pero luego esta !la realidad! Solo tenemos que pasar un rato en su compania para recordar sus extranas y molestas costumbres. Nos resulto muy comodo olvidar los rasgos de su caracter que nos sacaban de quicio. Pero eso no significa que los queramos menos. En la oscuridad de la Luna tienes que adaptarte a una realidad que no se ajusta a una fantasia. Pero ver a otra persona tal y como es en realidad tiene sus recompensas. Puede haber una conexion interesante.'
$ws.Rows.Item(12).AutoFit()

$ws.Range('B13').Value = 'PredicciA3n
                        para hoy                        15 de Agosto de 2023
Hay muchas cosas que nunca llegaremos a saber. Entre ellas estan las que creemos saber (pero que desconocemos). Hacemos conjeturas basandonos en informacion erronea y despues amontonamos capas de logica encima de ellas. Luego nos preguntamos por que las cosas no nos salen como habiamos planeado. Es un milagro que esto no ocurra mas a menudo. Estas reestructurando un plan para tu futuro. Es todo un reto. Pero si aplicas todo lo que has aprendido y la experiencia adquirida hasta ahora, construiras un camino a seguir con cimientos solidos que perduraran en el tiempo.'
$ws.Rows.Item(13).AutoFit()

